$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'309.57"
$ws.Range("E2").Value = "'-3.44%"
$ws.Range("G2").Value = "'11"
# Row 3
$ws.Range("D3").Value = "'50.63"
$ws.Range("E3").Value = "'3.59%"
$ws.Range("G3").Value = "'11"
# Row 4
$ws.Range("E4").Value = "'-1.49%"
$ws.Range("G4").Value = "'11"
# Row 5
$ws.Range("D5").Value = "'0.07772"
$ws.Range("E5").Value = "'-4.18%"
$ws.Range("G5").Value = "'11"
# Row 6
$ws.Range("D6").Value = "'4.497"
$ws.Range("E6").Value = "'-2.13%"
$ws.Range("G6").Value = "'11"
# Row 7
$ws.Range("D7").Value = "'1.342"
$ws.Range("E7").Value = "'11.11%"
$ws.Range("G7").Value = "'11"
# Row 8
$ws.Range("D8").Value = "'1.561"
$ws.Range("E8").Value = "'-5.23%"
$ws.Range("G8").Value = "'11"
# Row 9
$ws.Range("D9").Value = "'0.1212"
$ws.Range("E9").Value = "'-6.23%"
$ws.Range("G9").Value = "'11"
# Row 10
$ws.Range("D10").Value = "'0.1977"
$ws.Range("G10").Value = "'11"
# Row 11
$ws.Range("D11").Value = "'0.04790"
$ws.Range("E11").Value = "'4.53%"
$ws.Range("G11").Value = "'11"
# Row 12
$ws.Range("D12").Value = "'0.09412"
$ws.Range("E12").Value = "'-0.67%"
$ws.Range("G12").Value = "'11"
# Row 13
$ws.Range("E13").Value = "'-0.57%"
$ws.Range("G13").Value = "'11"
# Row 14
$ws.Range("D14").Value = "'0.001256"
$ws.Range("E14").Value = "'-5.45%"
$ws.Range("G14").Value = "'11"
# Row 15
$ws.Range("D15").Value = "'0.005785"
$ws.Range("E15").Value = "'-1.13%"
$ws.Range("G15").Value = "'11"
# Row 16
$ws.Range("E16").Value = "'2,016.20%"
$ws.Range("G16").Value = "'11"
# Row 17
$ws.Range("D17").Value = "'3.329"
$ws.Range("E17").Value = "'-0.36%"
$ws.Range("G17").Value = "'11"
# Row 18
$ws.Range("D18").Value = "'2.437"
$ws.Range("E18").Value = "'0.26%"
$ws.Range("G18").Value = "'11"
# Row 19
$ws.Range("E19").Value = "'1.65%"
$ws.Range("G19").Value = "'11"
# Row 20
$ws.Range("D20").Value = "'8.022"
$ws.Range("E20").Value = "'-0.80%"
$ws.Range("G20").Value = "'11"
# Row 21
$ws.Range("D21").Value = "'0.1363"
$ws.Range("E21").Value = "'-2.00%"
$ws.Range("G21").Value = "'11"
# Row 22
$ws.Range("D22").Value = "'0.3094"
$ws.Range("E22").Value = "'-1.00%"
$ws.Range("G22").Value = "'11"
# Row 23
$ws.Range("D23").Value = "'0.04173"
$ws.Range("E23").Value = "'-0.01%"
$ws.Range("G23").Value = "'11"
# Row 24
$ws.Range("D24").Value = "'0.001271"
$ws.Range("E24").Value = "'-2.61%"
$ws.Range("G24").Value = "'11"
# Row 25
$ws.Range("D25").Value = "'0.003949"
$ws.Range("E25").Value = "'-7.19%"
$ws.Range("G25").Value = "'11"
# Row 26
$ws.Range("D26").Value = "'0.0001349"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("G26").Value = "'11"
# Row 27
$ws.Range("G27").Value = "'11"
# Row 28
$ws.Range("G28").Value = "'11"
# Row 29
$ws.Range("G29").Value = "'11"
# Row 30
$ws.Range("G30").Value = "'11"
# Row 31
$ws.Range("G31").Value = "'11"
# Row 32
$ws.Range("G32").Value = "'11"
# Row 33
$ws.Range("G33").Value = "'11"
# Row 34
$ws.Range("G34").Value = "'11"
# Row 35
$ws.Range("G35").Value = "'11"
# Row 36
$ws.Range("G36").Value = "'11"
# Row 37
$ws.Range("G37").Value = "'11"
# Row 38
$ws.Range("D38").Value = "'0.02604"
$ws.Range("E38").Value = "'-4.01%"
$ws.Range("G38").Value = "'11"
# Row 39
$ws.Range("D39").Value = "'0.06119"
$ws.Range("E39").Value = "'6.67%"
$ws.Range("G39").Value = "'11"
# Row 40
$ws.Range("D40").Value = "'0.01099"
$ws.Range("E40").Value = "'74.50%"
$ws.Range("G40").Value = "'11"
# Row 41
$ws.Range("D41").Value = "'0.007926"
$ws.Range("E41").Value = "'2.75%"
$ws.Range("G41").Value = "'11"
# Row 42
$ws.Range("D42").Value = "'0.1421"
$ws.Range("E42").Value = "'-1.32%"
$ws.Range("G42").Value = "'11"
# Row 43
$ws.Range("D43").Value = "'0.008396"
$ws.Range("E43").Value = "'9.33%"
$ws.Range("G43").Value = "'11"
# Row 44
$ws.Range("D44").Value = "'0.008347"
$ws.Range("E44").Value = "'3.03%"
$ws.Range("G44").Value = "'11"
# Row 45
$ws.Range("D45").Value = "'0.3368"
$ws.Range("E45").Value = "'5.44%"
$ws.Range("G45").Value = "'11"
# Row 46
$ws.Range("D46").Value = "'0.00007204"
$ws.Range("E46").Value = "'3.19%"
$ws.Range("G46").Value = "'11"
# Row 47
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("G47").Value = "'11"
# Row 48
$ws.Range("D48").Value = "'0.002619"
$ws.Range("E48").Value = "'-34.54%"
$ws.Range("G48").Value = "'11"
# Row 49
$ws.Range("E49").Value = "'-17.28%"
$ws.Range("G49").Value = "'11"
# Row 50
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("G50").Value = "'11"
# Row 51
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.05%"
$ws.Range("G51").Value = "'11"
